$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.210958904109589
$ws.Range("C2").Value = 0.5123287671232877
$ws.Range("J2").Value = 0.0273972602739726
$ws.Range("P2").Value = 0.1561643835616438
$ws.Range("S2").Value = 0.09315068493150686
$ws.Range("B3").Value = 0.005208333333333333
$ws.Range("C3").Value = 0.01041666666666667
$ws.Range("J3").Value = 0.03125
$ws.Range("P3").Value = 0.71875
$ws.Range("S3").Value = 0.234375
$ws.Range("B6").Value = 0.06030150753768844
$ws.Range("D6").Value = 0.005025125628140704
$ws.Range("F6").Value = 0.04020100502512563
$ws.Range("J6").Value = 0.3015075376884422
$ws.Range("O6").Value = 0.03015075376884422
$ws.Range("Q6").Value = 0.1155778894472362
$ws.Range("R6").Value = 0.07537688442211055
$ws.Range("S6").Value = 0.3718592964824121
$ws.Range("B7").Value = 0.1221719457013575
$ws.Range("D7").Value = 0.03619909502262444
$ws.Range("F7").Value = 0.03167420814479638
$ws.Range("J7").Value = 0.1402714932126697
$ws.Range("O7").Value = 0.03167420814479638
$ws.Range("Q7").Value = 0.1764705882352941
$ws.Range("R7").Value = 0.05882352941176471
$ws.Range("S7").Value = 0.4027149321266968
$ws.Range("B8").Value = 0.1184834123222749
$ws.Range("D8").Value = 0.03317535545023697
$ws.Range("F8").Value = 0.04028436018957346
$ws.Range("J8").Value = 0.1066350710900474
$ws.Range("O8").Value = 0.01658767772511848
$ws.Range("Q8").Value = 0.1587677725118483
$ws.Range("R8").Value = 0.0924170616113744
$ws.Range("S8").Value = 0.433649289099526
$ws.Range("B9").Value = 0.131578947368421
$ws.Range("D9").Value = 0.01973684210526316
$ws.Range("F9").Value = 0.06578947368421052
$ws.Range("J9").Value = 0.09210526315789473
$ws.Range("O9").Value = 0.02631578947368421
$ws.Range("Q9").Value = 0.1907894736842105
$ws.Range("R9").Value = 0.07236842105263158
$ws.Range("S9").Value = 0.4013157894736842
$ws.Range("B10").Value = 0.1304347826086956
$ws.Range("D10").Value = 0.03610906411201179
$ws.Range("F10").Value = 0.06263817243920412
$ws.Range("J10").Value = 0.1149594694178335
$ws.Range("O10").Value = 0.0103168754605748
$ws.Range("Q10").Value = 0.2063375092114959
$ws.Range("R10").Value = 0.05305821665438467
$ws.Range("S10").Value = 0.3861459100957996
$ws.Range("G11").Value = 0.1289398280802292
$ws.Range("J11").Value = 0.09742120343839542
$ws.Range("K11").Value = 0.1833810888252149
$ws.Range("L11").Value = 0.5845272206303725
$ws.Range("S11").Value = 0.005730659025787965
$ws.Range("G12").Value = 0.6976744186046512
$ws.Range("J12").Value = 0.2465116279069768
$ws.Range("K12").Value = 0.009302325581395349
$ws.Range("L12").Value = 0.0186046511627907
$ws.Range("S12").Value = 0.02790697674418605
$ws.Range("F15").Value = 0.02252252252252252
$ws.Range("H15").Value = 0.1486486486486487
$ws.Range("I15").Value = 0.04504504504504504
$ws.Range("J15").Value = 0.3828828828828829
$ws.Range("K15").Value = 0.09009009009009009
$ws.Range("M15").Value = 0.01351351351351351
$ws.Range("O15").Value = 0.06756756756756757
$ws.Range("S15").Value = 0.2297297297297297
$ws.Range("F16").Value = 0.01282051282051282
$ws.Range("H16").Value = 0.1495726495726496
$ws.Range("I16").Value = 0.06837606837606838
$ws.Range("J16").Value = 0.405982905982906
$ws.Range("K16").Value = 0.1239316239316239
$ws.Range("M16").Value = 0.04273504273504274
$ws.Range("O16").Value = 0.04273504273504274
$ws.Range("S16").Value = 0.1538461538461539
$ws.Range("F17").Value = 0.01822323462414579
$ws.Range("H17").Value = 0.1435079726651481
$ws.Range("I17").Value = 0.06150341685649203
$ws.Range("J17").Value = 0.4646924829157175
$ws.Range("K17").Value = 0.1161731207289294
$ws.Range("M17").Value = 0.01594533029612756
$ws.Range("O17").Value = 0.05922551252847381
$ws.Range("S17").Value = 0.1207289293849658
$ws.Range("F18").Value = 0.03947368421052631
$ws.Range("H18").Value = 0.1578947368421053
$ws.Range("I18").Value = 0.06578947368421052
$ws.Range("J18").Value = 0.4407894736842105
$ws.Range("K18").Value = 0.08552631578947369
$ws.Range("M18").Value = 0.03289473684210526
$ws.Range("O18").Value = 0.07236842105263158
$ws.Range("S18").Value = 0.1052631578947368
$ws.Range("F19").Value = 0.01844734819369715
$ws.Range("H19").Value = 0.2044581091468101
$ws.Range("I19").Value = 0.06917755572636433
$ws.Range("J19").Value = 0.3858570330514988
$ws.Range("K19").Value = 0.1245196003074558
$ws.Range("M19").Value = 0.02613374327440431
$ws.Range("O19").Value = 0.07302075326671791
$ws.Range("S19").Value = 0.0983858570330515
